# Azerbaijan Premier League base update (17-02-2024 11:11)
#
# The upstream feed re-fetched the league's fixture list. Rows 82-85 held
# four already-finished matches that got re-ordered/re-scraped (same
# match ids, odds refreshed), so their B..AC payload rotates down by one
# slot (old row 85 -> row 82, old row 82 -> row 83, old row 83 -> row 84,
# old row 84 -> row 85) while the running id in column A and the shared
# Div / Date columns (A, C, D, E) stay put. A brand-new, not-yet-played
# fixture is appended as row 200.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Rows 82-85: rotate the match payload (everything except id/Div/Date)
# down by one row, wrapping the last row's data back up to the first.
# ---------------------------------------------------------------------

# New row 82 <= old row 85's data
$ws.Range("B82").Value = 5573342
$ws.Range("F82").Value = "PFK Turan Tovuz"
$ws.Range("G82").Value = "Sabail FC"
$ws.Range("H82").Value = 2
$ws.Range("I82").Value = 2
$ws.Range("J82").Value = "D"
$ws.Range("K82").Value = 2.6
$ws.Range("L82").Value = 3
$ws.Range("M82").Value = 2.6
$ws.Range("N82").Value = 2.8
$ws.Range("O82").Value = 2.875
$ws.Range("P82").Value = 2.5
$ws.Range("Q82").Value = 0
$ws.Range("R82").Value = 2.05
$ws.Range("S82").Value = 1.75
$ws.Range("T82").Value = 2.25
$ws.Range("U82").Value = 1.875
$ws.Range("V82").Value = 1.925
$ws.Range("W82").Value = -1
$ws.Range("X82").Value = 1.875
$ws.Range("Y82").Value = -1
$ws.Range("Z82").Value = 0
$ws.Range("AA82").Value = -0
$ws.Range("AB82").Value = 0.875
$ws.Range("AC82").Value = -1

# New row 83 <= old row 82's data
$ws.Range("B83").Value = 5579144
$ws.Range("F83").Value = "Sabah"
$ws.Range("G83").Value = "Zira IK"
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = "D"
$ws.Range("K83").Value = 1.45
$ws.Range("L83").Value = 4.2
$ws.Range("M83").Value = 5.5
$ws.Range("N83").Value = 1.5
$ws.Range("O83").Value = 4
$ws.Range("P83").Value = 5.25
$ws.Range("Q83").Value = -1
$ws.Range("R83").Value = 1.85
$ws.Range("S83").Value = 1.95
$ws.Range("T83").Value = 2.5
$ws.Range("U83").Value = 1.8
$ws.Range("V83").Value = 2
$ws.Range("W83").Value = -1
$ws.Range("X83").Value = 3
$ws.Range("Y83").Value = -1
$ws.Range("Z83").Value = -1
$ws.Range("AA83").Value = 0.95
$ws.Range("AB83").Value = -1
$ws.Range("AC83").Value = 1

# New row 84 <= old row 83's data
$ws.Range("B84").Value = 5574442
$ws.Range("F84").Value = "FK Qarabag"
$ws.Range("G84").Value = "FK Sumqayit"
$ws.Range("H84").Value = 1
$ws.Range("I84").Value = 2
$ws.Range("J84").Value = "A"
$ws.Range("K84").Value = 1.125
$ws.Range("L84").Value = 7.5
$ws.Range("M84").Value = 15
$ws.Range("N84").Value = 1.2
$ws.Range("O84").Value = 6
$ws.Range("P84").Value = 11
$ws.Range("Q84").Value = -2.25
$ws.Range("R84").Value = 1.975
$ws.Range("S84").Value = 1.825
$ws.Range("T84").Value = 3.5
$ws.Range("U84").Value = 1.825
$ws.Range("V84").Value = 1.975
$ws.Range("W84").Value = -1
$ws.Range("X84").Value = -1
$ws.Range("Y84").Value = 10
$ws.Range("Z84").Value = -1
$ws.Range("AA84").Value = 0.825
$ws.Range("AB84").Value = -1
$ws.Range("AC84").Value = 0.9750000000000001

# New row 85 <= old row 84's data
$ws.Range("B85").Value = 5573343
$ws.Range("F85").Value = "Shamakhi FK"
$ws.Range("G85").Value = "FK Gabala"
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = "D"
$ws.Range("K85").Value = 3.5
$ws.Range("L85").Value = 3.1
$ws.Range("M85").Value = 2
$ws.Range("N85").Value = 3.3
$ws.Range("O85").Value = 3.2
$ws.Range("P85").Value = 2.05
$ws.Range("Q85").Value = 0.25
$ws.Range("R85").Value = 2
$ws.Range("S85").Value = 1.8
$ws.Range("T85").Value = 2.5
$ws.Range("U85").Value = 1.975
$ws.Range("V85").Value = 1.825
$ws.Range("W85").Value = -1
$ws.Range("X85").Value = 2.2
$ws.Range("Y85").Value = -1
$ws.Range("Z85").Value = 0.5
$ws.Range("AA85").Value = -0.5
$ws.Range("AB85").Value = -1
$ws.Range("AC85").Value = 0.825

# ---------------------------------------------------------------------
# Append row 200: a brand-new (not-yet-played) fixture.
# Copy the format of the last existing data row (199) first so the
# bold/bordered id style and the custom date format follow along, then
# fill in the values. FTHG/FTAG/FTR and the PL_* closing columns are
# left blank since the match hasn't been played yet.
# ---------------------------------------------------------------------

$ws.Range("A199:AC199").Copy()
$ws.Range("A200:AC200").PasteSpecial(-4122)
$ws.Range("H200").ClearContents()
$ws.Range("I200").ClearContents()
$ws.Range("J200").ClearContents()
$ws.Range("AB200").ClearContents()
$ws.Range("AC200").ClearContents()

$ws.Range("A200").Value = 198
$ws.Range("B200").Value = 7011602
$ws.Range("C200").Value = "Azerbaijan Premier League"
$ws.Range("D200").Value = "Azerbaijan Premier League"
$ws.Range("E200").Value = 45339.47916666666
$ws.Range("F200").Value = "Sabail FC"
$ws.Range("G200").Value = "Araz FK"
$ws.Range("K200").Value = 1.8
$ws.Range("L200").Value = 3
$ws.Range("M200").Value = 4.5
$ws.Range("N200").Value = 2.3
$ws.Range("O200").Value = 3.1
$ws.Range("P200").Value = 2.875
$ws.Range("Q200").Value = -0.25
$ws.Range("R200").Value = 2.025
$ws.Range("S200").Value = 1.775
$ws.Range("T200").Value = 2.5
$ws.Range("U200").Value = 1.95
$ws.Range("V200").Value = 1.85
$ws.Range("W200").Value = 0
$ws.Range("X200").Value = 0
$ws.Range("Y200").Value = 0
$ws.Range("Z200").Value = 0
$ws.Range("AA200").Value = 0
